$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the standard data-row style (s=3), used to restore
# formatting on cells where Excel auto-applies a Percentage number format
# after a plain "NN%" value is typed in (matches native Excel behaviour).
$fmtSource = $ws.Range("C2")

$ws.Range('E2').Value = '2026-02-14 19:48:52'
$ws.Range('I2').Value = '34.6 mm'
$ws.Range('N2').Value = '-3.1 °C 19:19 TU'
$ws.Range('O2').Value = '-1.1 °C'
$ws.Range('E3').Value = '2026-02-14 19:48:55'
$ws.Range('I3').Value = '14.7 mm'
$ws.Range('N3').Value = '-6.7 °C 19:08 TU'
$ws.Range('O3').Value = '-5.1 °C'
$ws.Range('E4').Value = '2026-02-14 19:48:58'
$ws.Range('J4').Value = '996.5 hPa'
$ws.Range('N4').Value = '6.0 °C 19:15 TU'
$ws.Range('O4').Value = '10.9 °C'
$ws.Range('E5').Value = '2026-02-14 19:49:00'
$ws.Range('I5').Value = '21.3 mm'
$ws.Range('N5').Value = '-6.6 °C 19:28 TU'
$ws.Range('E6').Value = '2026-02-14 19:49:02'
$ws.Range('H6').NumberFormat = "@"
$ws.Range('H6').Value = '77%'
$fmtSource.Copy()
$ws.Range('H6').PasteSpecial(-4122)
$ws.Range('J6').Value = '996.5 hPa'
$ws.Range('O6').Value = '10.2 °C'
$ws.Range('E7').Value = '2026-02-14 19:49:05'
$ws.Range('H7').NumberFormat = "@"
$ws.Range('H7').Value = '53%'
$fmtSource.Copy()
$ws.Range('H7').PasteSpecial(-4122)
$ws.Range('J7').Value = '996.7 hPa'
$ws.Range('E8').Value = '2026-02-14 19:49:08'
$ws.Range('J8').Value = '996.5 hPa'
$ws.Range('O8').Value = '9.8 °C'
$ws.Range('E9').Value = '2026-02-14 19:49:11'
$ws.Range('N9').Value = '10.5 °C 19:07 TU'
$ws.Range('E10').Value = '2026-02-14 19:49:14'
$ws.Range('H10').NumberFormat = "@"
$ws.Range('H10').Value = '78%'
$fmtSource.Copy()
$ws.Range('H10').PasteSpecial(-4122)
$ws.Range('E11').Value = '2026-02-14 19:49:16'
$ws.Range('E12').Value = '2026-02-14 19:49:19'
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H12').Value = '55%'
$fmtSource.Copy()
$ws.Range('H12').PasteSpecial(-4122)
$ws.Range('O12').Value = '12.2 °C'
$ws.Range('E13').Value = '2026-02-14 19:49:22'
$ws.Range('H13').NumberFormat = "@"
$ws.Range('H13').Value = '67%'
$fmtSource.Copy()
$ws.Range('H13').PasteSpecial(-4122)
$ws.Range('J13').Value = '999.3 hPa'
$ws.Range('O13').Value = '4.0 °C'
$ws.Range('E14').Value = '2026-02-14 19:49:24'
$ws.Range('K14').Value = '13.9 MJ/m2'
$ws.Range('E15').Value = '2026-02-14 19:49:27'
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H15').Value = '57%'
$fmtSource.Copy()
$ws.Range('H15').PasteSpecial(-4122)
$ws.Range('E16').Value = '2026-02-14 19:49:30'
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H16').Value = '76%'
$fmtSource.Copy()
$ws.Range('H16').PasteSpecial(-4122)
$ws.Range('N16').Value = '-8.7 °C 19:17 TU'
$ws.Range('E17').Value = '2026-02-14 19:49:33'
$ws.Range('H17').NumberFormat = "@"
$ws.Range('H17').Value = '68%'
$fmtSource.Copy()
$ws.Range('H17').PasteSpecial(-4122)
$ws.Range('E18').Value = '2026-02-14 19:49:35'
$ws.Range('H18').NumberFormat = "@"
$ws.Range('H18').Value = '76%'
$fmtSource.Copy()
$ws.Range('H18').PasteSpecial(-4122)
$ws.Range('J18').Value = '996.7 hPa'
$ws.Range('O18').Value = '10.7 °C'
$ws.Range('E19').Value = '2026-02-14 19:49:38'
$ws.Range('H19').NumberFormat = "@"
$ws.Range('H19').Value = '79%'
$fmtSource.Copy()
$ws.Range('H19').PasteSpecial(-4122)
$ws.Range('O19').Value = '5.7 °C'
$ws.Range('E20').Value = '2026-02-14 19:49:41'
$ws.Range('I20').Value = '4.0 mm'
$ws.Range('N20').Value = '-7.4 °C 19:24 TU'
$ws.Range('E21').Value = '2026-02-14 19:49:44'
$ws.Range('J21').Value = '999.1 hPa'
$ws.Range('E22').Value = '2026-02-14 19:49:46'
$ws.Range('N22').Value = '-8.7 °C 19:19 TU'
$ws.Range('E23').Value = '2026-02-14 19:49:49'
$ws.Range('I23').Value = '38.2 mm'
$ws.Range('N23').Value = '-8.3 °C 19:15 TU'
$ws.Range('E24').Value = '2026-02-14 19:49:52'
$ws.Range('H24').NumberFormat = "@"
$ws.Range('H24').Value = '66%'
$fmtSource.Copy()
$ws.Range('H24').PasteSpecial(-4122)
$ws.Range('J24').Value = '1000.7 hPa'
$ws.Range('E25').Value = '2026-02-14 19:49:55'
$ws.Range('I25').Value = '13.7 mm'
$ws.Range('N25').Value = '-7.0 °C 19:13 TU'
$ws.Range('O25').Value = '-4.6 °C'
$ws.Range('E26').Value = '2026-02-14 19:49:58'
$ws.Range('E27').Value = '2026-02-14 19:50:00'
$ws.Range('N27').Value = '-5.9 °C 19:29 TU'
$ws.Range('E28').Value = '2026-02-14 19:50:03'
$ws.Range('H28').NumberFormat = "@"
$ws.Range('H28').Value = '69%'
$fmtSource.Copy()
$ws.Range('H28').PasteSpecial(-4122)
$ws.Range('J28').Value = '996.5 hPa'
$ws.Range('O28').Value = '9.2 °C'
$ws.Range('E29').Value = '2026-02-14 19:50:06'
$ws.Range('H29').NumberFormat = "@"
$ws.Range('H29').Value = '64%'
$fmtSource.Copy()
$ws.Range('H29').PasteSpecial(-4122)
$ws.Range('E30').Value = '2026-02-14 19:50:09'
$ws.Range('H30').NumberFormat = "@"
$ws.Range('H30').Value = '57%'
$fmtSource.Copy()
$ws.Range('H30').PasteSpecial(-4122)
$ws.Range('J30').Value = '996.4 hPa'
$ws.Range('O30').Value = '11.6 °C'
$ws.Range('E31').Value = '2026-02-14 19:50:12'
$ws.Range('J31').Value = '995.5 hPa'
$ws.Range('K31').Value = '6.8 MJ/m2'
$ws.Range('E32').Value = '2026-02-14 19:50:15'
$ws.Range('E33').Value = '2026-02-14 19:50:18'
$ws.Range('J33').Value = '998.5 hPa'
$ws.Range('L33').Value = '54.7 km/h - 323º 19:22 TU'
$ws.Range('E34').Value = '2026-02-14 19:50:20'
$ws.Range('N34').Value = '-4.5 °C 19:20 TU'
$ws.Range('O34').Value = '-2.2 °C'
$ws.Range('E35').Value = '2026-02-14 19:50:23'
$ws.Range('J35').Value = '1003.4 hPa'
$ws.Range('E36').Value = '2026-02-14 19:50:26'
$ws.Range('J36').Value = '997.2 hPa'
$ws.Range('L36').Value = '92.5 km/h - 320º 19:21 TU'
$ws.Range('N36').Value = '10.9 °C 19:17 TU'
$ws.Range('E37').Value = '2026-02-14 19:50:29'
$ws.Range('H37').NumberFormat = "@"
$ws.Range('H37').Value = '66%'
$fmtSource.Copy()
$ws.Range('H37').PasteSpecial(-4122)
$ws.Range('J37').Value = '997.4 hPa'
$ws.Range('L37').Value = '61.2 km/h - 41º 19:04 TU'
$ws.Range('E38').Value = '2026-02-14 19:50:31'
$ws.Range('E39').Value = '2026-02-14 19:50:34'
$ws.Range('I39').Value = '12.3 mm'
$ws.Range('N39').Value = '-8.4 °C 19:28 TU'
$ws.Range('E40').Value = '2026-02-14 19:50:37'
$ws.Range('H40').NumberFormat = "@"
$ws.Range('H40').Value = '66%'
$fmtSource.Copy()
$ws.Range('H40').PasteSpecial(-4122)
$ws.Range('J40').Value = '999.7 hPa'
$ws.Range('O40').Value = '7.1 °C'
$ws.Range('E41').Value = '2026-02-14 19:50:40'
$ws.Range('J41').Value = '998.5 hPa'
$ws.Range('E42').Value = '2026-02-14 19:50:43'
$ws.Range('E43').Value = '2026-02-14 19:50:45'
$ws.Range('E44').Value = '2026-02-14 19:50:47'
$ws.Range('G44').Value = '269 cm'
$ws.Range('I44').Value = '36.4 mm'
$ws.Range('N44').Value = '-6.8 °C 19:29 TU'
$ws.Range('E45').Value = '2026-02-14 19:50:50'
$ws.Range('H45').NumberFormat = "@"
$ws.Range('H45').Value = '82%'
$fmtSource.Copy()
$ws.Range('H45').PasteSpecial(-4122)
$ws.Range('J45').Value = '1005.6 hPa'
$ws.Range('E46').Value = '2026-02-14 19:50:53'
$ws.Range('J46').Value = '1001.6 hPa'

$excel.CutCopyMode = $false
